$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D (Price) need Text number format so Excel does not
# reinterpret values like "1.00", "11.40" or "0.0000141" as numbers and
# strip formatting/precision.
$priceCells = @("D2","D3","D5","D6","D7","D9","D11","D12","D13","D14","D15","D16","D17","D18","D19","D21","D22","D24","D25","D27","D28","D29","D30","D31","D32","D33","D36","D39","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "60.355.09"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").Value = "2.600.44"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "585.19"
$ws.Range("E5").Value = "  +5.18%  "
$ws.Range("D6").Value = "142.33"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("D9").Value = "2.608.34"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("E10").Value = "  -3.10%  "
$ws.Range("D11").Value = "0.105"
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").Value = "0.155"
$ws.Range("E12").Value = "  -3.99%  "
$ws.Range("D13").Value = "0.371"
$ws.Range("E13").Value = "  +4.44%  "
$ws.Range("D14").Value = "3.073.68"
$ws.Range("E14").Value = "  +1.32%  "
$ws.Range("D15").Value = "24.63"
$ws.Range("E15").Value = "  +5.36%  "
$ws.Range("D16").Value = "60.379.77"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").Value = "0.0000141"
$ws.Range("E17").Value = "  +2.21%  "
$ws.Range("D18").Value = "2.614.38"
$ws.Range("E18").Value = "  +1.23%  "
$ws.Range("D19").Value = "11.40"
$ws.Range("E19").Value = "  +9.48%  "
$ws.Range("E20").Value = "  +1.90%  "
$ws.Range("D21").Value = "347.53"
$ws.Range("E21").Value = "  +2.11%  "
$ws.Range("D22").Value = "6.91"
$ws.Range("E22").Value = "  +6.02%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "0.522"
$ws.Range("E24").Value = "  +8.91%  "
$ws.Range("D25").Value = "62.98"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("E26").Value = "  +0.42%  "
$ws.Range("D27").Value = "0.159"
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("D28").Value = "8.02"
$ws.Range("E28").Value = "  +7.21%  "
$ws.Range("D29").Value = "0.0₃0793"
$ws.Range("E29").Value = "  +2.46%  "
$ws.Range("D30").Value = "1.86"
$ws.Range("E30").Value = "  +10.38%  "
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").Value = "6.38"
$ws.Range("E31").Value = "  +2.56%  "
$ws.Range("B32").Value = "USDe"
$ws.Range("C32").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D32").Value = "0.998"
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("D33").Value = "163.27"
$ws.Range("E33").Value = "  +3.28%  "
$ws.Range("E34").Value = "  +1.54%  "
$ws.Range("E35").Value = "  +2.15%  "
$ws.Range("D36").Value = "0.980"
$ws.Range("E36").Value = "  +7.13%  "
$ws.Range("E37").Value = "  +4.60%  "
$ws.Range("E38").Value = "  +8.45%  "
$ws.Range("D39").Value = "37.91"
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("E40").Value = "  +5.75%  "
$ws.Range("D41").Value = "310.86"
$ws.Range("E41").Value = "  +6.77%  "
$ws.Range("D42").Value = "0.840"
$ws.Range("E42").Value = "  -1.32%  "
$ws.Range("D43").Value = "135.58"
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "0.0994"
$ws.Range("E44").Value = "  +2.06%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "19.87"
$ws.Range("E46").Value = "  +4.38%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "0.604"
$ws.Range("E47").Value = "  +0.73%  "
$ws.Range("D48").Value = "0.0549"
$ws.Range("E48").Value = "  +2.93%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "20.17"
$ws.Range("E49").Value = "  +7.08%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "4.99"
$ws.Range("E50").Value = "  +4.89%  "
$ws.Range("D51").Value = "0.0241"
$ws.Range("E51").Value = "  +2.47%  "
